$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.072.38"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.379.25"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'302.02"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'96.92"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'0.505"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").Value = "'34.32"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").Value = "'18.25"
$ws.Range("E13").Value = "  -5.00%  "
$ws.Range("D14").Value = "'6.81"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "2.748.33"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "2.385.90"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "'0.805"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").Value = "43.026.19"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").Value = "'12.20"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "'6.32"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").Value = "0.0₃0888"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "'68.14"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'235.52"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "'2.25"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.44"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'24.96"
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").Value = "'2.36"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'9.26"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "'31.52"
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "'5.08"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("E33").Value = "  +5.13%  "
$ws.Range("D34").Value = "'17.57"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.106"
$ws.Range("E35").Value = "  +5.52%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.87"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("D37").Value = "'4.38"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.81"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'22.65"
$ws.Range("E40").Value = "  +10.35%  "
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").Value = "'115.14"
$ws.Range("E42").Value = "  -30.84%  "
$ws.Range("D43").Value = "1.953.23"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'0.0281"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").Value = "'2.75"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("D47").Value = "'9.17"
$ws.Range("E47").Value = "  -11.42%  "
$ws.Range("D48").Value = "'1.53"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").Value = "'52.43"
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("D50").Value = "'72.41"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  +0.84%  "
